# Adds 33 new arrival rows (Friday, Jan 13 schedule) to the "Main Data" sheet,
# appending after the existing last row (row 395 / NUMBER 394).
# Columns: A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=FROM, F=SHORT, G=AIRLINE,
#          H=MODEL, I=AIRCRAFT ID, J=STATUS, K=(blank), L=DIFFERENCE, M=(blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    [PSCustomObject]@{ A=395.0; B='Friday, Jan 13'; C='7:41 AM'; D='LO9001'; E='Warsaw'; F='(WAW)'; G='LOT '; H='B738'; I='(SP-LWC)'; J='8:10 AM'; L='0 hours, 29 minutes' },
    [PSCustomObject]@{ A=396.0; B='Friday, Jan 13'; C='8:10 AM'; D='LO3903'; E='Warsaw'; F='(WAW)'; G='LOT '; H='E75S'; I='(SP-LIK)'; J='8:16 AM'; L='0 hours, 6 minutes' },
    [PSCustomObject]@{ A=397.0; B='Friday, Jan 13'; C='8:15 AM'; D='LG5741'; E='Luxembourg'; F='(LUX)'; G='Luxair '; H='DH8D'; I='(LX-LQB)'; J='7:54 AM'; L='0 hours, -21 minutes' },
    [PSCustomObject]@{ A=398.0; B='Friday, Jan 13'; C='8:45 AM'; D='LX1370'; E='Zurich'; F='(ZRH)'; G='Air Baltic '; H='BCS3'; I='(YL-ABG)'; J='8:36 AM'; L='0 hours, -9 minutes' },
    [PSCustomObject]@{ A=399.0; B='Friday, Jan 13'; C='9:10 AM'; D='FR2228'; E='Valencia'; F='(VLC)'; G='Ryanair '; H='B738'; I='(EI-DPZ)'; J='8:41 AM'; L='0 hours, -29 minutes' },
    [PSCustomObject]@{ A=400.0; B='Friday, Jan 13'; C='9:30 AM'; D='FR3053'; E='Barcelona'; F='(BCN)'; G='Ryanair '; H='B738'; I='(EI-ENX)'; J='9:08 AM'; L='0 hours, -22 minutes' },
    [PSCustomObject]@{ A=401.0; B='Friday, Jan 13'; C='9:35 AM'; D='FR6293'; E='London'; F='(STN)'; G='Ryanair '; H='B38M'; I='(EI-HMW)'; J='9:20 AM'; L='0 hours, -15 minutes' },
    [PSCustomObject]@{ A=402.0; B='Friday, Jan 13'; C='10:00 AM'; D='U26939'; E='Edinburgh'; F='(EDI)'; G='easyJet '; H='A319'; I='(G-EZBF)'; J='9:42 AM'; L='0 hours, -18 minutes' },
    [PSCustomObject]@{ A=403.0; B='Friday, Jan 13'; C='10:10 AM'; D='FR5118'; E='Dublin'; F='(DUB)'; G='Ryanair '; H='B38M'; I='(EI-HMV)'; J='9:39 AM'; L='0 hours, -31 minutes' },
    [PSCustomObject]@{ A=404.0; B='Friday, Jan 13'; C='10:15 AM'; D='FR3986'; E='Bari'; F='(BRI)'; G='Ryanair '; H='B738'; I='(SP-RSA)'; J='10:02 AM'; L='0 hours, -13 minutes' },
    [PSCustomObject]@{ A=405.0; B='Friday, Jan 13'; C='10:23 AM'; D='UNKNOWN'; E='Nuremberg'; F='(NUE)'; G='ADAC Luftrettung '; H='J328'; I='(D-BADC)'; J='10:15 AM'; L='0 hours, -8 minutes' },
    [PSCustomObject]@{ A=406.0; B='Friday, Jan 13'; C='10:35 AM'; D='U26275'; E='Bristol'; F='(BRS)'; G='easyJet '; H='A320'; I='(G-EZUF)'; J='10:19 AM'; L='0 hours, -16 minutes' },
    [PSCustomObject]@{ A=407.0; B='Friday, Jan 13'; C='10:40 AM'; D='FR6225'; E='Podgorica'; F='(TGD)'; G='Buzz '; H='B38M'; I='(SP-RZA)'; J='10:37 AM'; L='0 hours, -3 minutes' },
    [PSCustomObject]@{ A=408.0; B='Friday, Jan 13'; C='10:40 AM'; D='LS353'; E='Newcastle'; F='(NCL)'; G='Jet2 '; H='B738'; I='(G-JZHM)'; J='10:25 AM'; L='0 hours, -15 minutes' },
    [PSCustomObject]@{ A=409.0; B='Friday, Jan 13'; C='10:55 AM'; D='AY1163'; E='Helsinki'; F='(HEL)'; G='Finnair '; H='E190'; I='(OH-LKF)'; J='10:41 AM'; L='0 hours, -14 minutes' },
    [PSCustomObject]@{ A=410.0; B='Friday, Jan 13'; C='11:00 AM'; D='FR6211'; E='Paris'; F='(BVA)'; G='Buzz '; H='B38M'; I='(SP-RZF)'; J='10:54 AM'; L='0 hours, -6 minutes' },
    [PSCustomObject]@{ A=411.0; B='Friday, Jan 13'; C='11:10 AM'; D='FZ1787'; E='Dubai'; F='(DXB)'; G='flydubai '; H='B38M'; I='(A6-FML)'; J='10:45 AM'; L='0 hours, -25 minutes' },
    [PSCustomObject]@{ A=412.0; B='Friday, Jan 13'; C='11:30 AM'; D='FR4934'; E='Brussels'; F='(CRL)'; G='Buzz '; H='B38M'; I='(SP-RZH)'; J='11:57 AM'; L='0 hours, 27 minutes' },
    [PSCustomObject]@{ A=413.0; B='Friday, Jan 13'; C='11:30 AM'; D='LO3907'; E='Warsaw'; F='(WAW)'; G='LOT (Grzeski Livery) '; H='E195'; I='(SP-LNB)'; J='11:20 AM'; L='0 hours, -10 minutes' },
    [PSCustomObject]@{ A=414.0; B='Friday, Jan 13'; C='11:30 AM'; D='W65034'; E='Oslo'; F='(OSL)'; G='Wizz Air '; H='A21N'; I='(HA-LVH)'; J='11:32 AM'; L='0 hours, 2 minutes' },
    [PSCustomObject]@{ A=415.0; B='Friday, Jan 13'; C='11:45 AM'; D='W65014'; E='Leeds'; F='(LBA)'; G='Wizz Air '; H='A21N'; I='(HA-LVO)'; J='12:11 PM'; L='0 hours, 26 minutes' },
    [PSCustomObject]@{ A=416.0; B='Friday, Jan 13'; C='12:05 PM'; D='DY1040'; E='Oslo'; F='(OSL)'; G='Norwegian '; H='B738'; I='(LN-NII)'; J='11:54 AM'; L='0 hours, -11 minutes' },
    [PSCustomObject]@{ A=417.0; B='Friday, Jan 13'; C='12:05 PM'; D='LH1366'; E='Frankfurt'; F='(FRA)'; G='Lufthansa '; H='A320'; I='(D-AIZY)'; J='12:54 PM'; L='0 hours, 49 minutes' },
    [PSCustomObject]@{ A=418.0; B='Friday, Jan 13'; C='12:10 PM'; D='FR6233'; E='Palermo'; F='(PMO)'; G='Ryanair '; H='B38M'; I='(SP-RZL)'; J='12:13 PM'; L='0 hours, 3 minutes' },
    [PSCustomObject]@{ A=419.0; B='Friday, Jan 13'; C='12:10 PM'; D='FR6249'; E='Manchester'; F='(MAN)'; G='Ryanair '; H='B738'; I='(SP-RKU)'; J='12:44 PM'; L='0 hours, 34 minutes' },
    [PSCustomObject]@{ A=420.0; B='Friday, Jan 13'; C='12:35 PM'; D='FR1813'; E='London'; F='(LTN)'; G='Buzz '; H='B38M'; I='(SP-RZB)'; J='12:21 PM'; L='0 hours, -14 minutes' },
    [PSCustomObject]@{ A=421.0; B='Friday, Jan 13'; C='12:35 PM'; D='FR4204'; E='Glasgow'; F='(GLA)'; G='Buzz '; H='B38M'; I='(SP-RZD)'; J='12:17 PM'; L='0 hours, -18 minutes' },
    [PSCustomObject]@{ A=422.0; B='Friday, Jan 13'; C='12:45 PM'; D='LH1620'; E='Munich'; F='(MUC)'; G='Lufthansa '; H='A320'; I='(D-AIWA)'; J='12:48 PM'; L='0 hours, 3 minutes' },
    [PSCustomObject]@{ A=423.0; B='Friday, Jan 13'; C='12:50 PM'; D='FR7954'; E='Prague'; F='(PRG)'; G='Ryanair '; H='B738'; I='(SP-RSH)'; J='12:38 PM'; L='0 hours, -12 minutes' },
    [PSCustomObject]@{ A=424.0; B='Friday, Jan 13'; C='1:10 PM'; D='W65048'; E='Barcelona'; F='(BCN)'; G='Wizz Air '; H='A21N'; I='(HA-LZI)'; J='12:52 PM'; L='0 hours, -18 minutes' },
    [PSCustomObject]@{ A=425.0; B='Friday, Jan 13'; C='1:50 PM'; D='OS597'; E='Vienna'; F='(VIE)'; G='Austrian Airlines '; H='E195'; I='(OE-LWA)'; J='1:33 PM'; L='0 hours, -17 minutes' },
    [PSCustomObject]@{ A=426.0; B='Friday, Jan 13'; C='1:55 PM'; D='FR9279'; E='Agadir'; F='(AGA)'; G='Ryanair '; H='B738'; I='(EI-EVA)'; J='1:41 PM'; L='0 hours, -14 minutes' },
    [PSCustomObject]@{ A=427.0; B='Friday, Jan 13'; C='2:20 PM'; D='LO3905'; E='Warsaw'; F='(WAW)'; G='LOT '; H='E195'; I='(SP-LNP)'; J='2:15 PM'; L='0 hours, -5 minutes' }
)

$startRow = 396
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 12).Value = $row.L
}
